$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- Insert a new row at 86, shifting existing rows (86-201) down to (87-202) ---
$ws.Rows("86:86").Insert()

# --- Grow the table (Table1) to include the freshly inserted row ---
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K202"))

# --- Restore the calculated column formula on the new row (EARNED  column) ---
$ws.Range("G86").Formula = '=IF(ISBLANK([@EARNED]),"",[@EARNED])'

# --- Populate the new leave-card entry (row 86) ---
$ws.Range("B86").Value = "A(1-0-0)"
$ws.Range("D86").Value = 1
$ws.Range("K86").Value = DateSerial(2023, 3, 7)

# --- Correct the particulars / date on row 84 ---
$ws.Range("B84").Value = "SP(1-0-0)"
$ws.Range("K84").Value = DateSerial(2023, 2, 20)
